$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 4 (pushes existing rows 4..29 down to 6..31)
$ws.Range("A4:A5").EntireRow.Insert()

# Row 4: new "Holden" entry
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Holden"
$ws.Range("C4").Value = 1.001634263994456
$ws.Range("D4").Value = 0.9934629611357327
$ws.Range("E4").Value = 1.001634263994456
$ws.Range("F4").Value = 0.9934629611357327
$ws.Range("G4").Value = 1.001634263994456
$ws.Range("H4").Value = 1.001634263994456
$ws.Range("I4").Value = 1.004358026082893
$ws.Range("J4").Value = 0.9964046282316398
$ws.Range("K4").Value = 1.001634263994456
$ws.Range("L4").Value = 1.001634263994456
$ws.Range("M4").Value = 0.9975486125650943
$ws.Range("N4").Value = 0.9975486125650943
$ws.Range("O4").Value = 0.9971672844539428
$ws.Range("P4").Value = 0.9989104963748815
$ws.Range("Q4").Value = 0.9989104963748815
$ws.Range("R4").Value = 0.999591438279775
$ws.Range("S4").Value = 0.999591438279775
$ws.Range("T4").Value = 0.999854734572272

# Row 5: new "Rizzie Spiral" entry
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "Rizzie Spiral"
$ws.Range("C5").Value = 1.001126565367507
$ws.Range("D5").Value = 0.995493715611908
$ws.Range("E5").Value = 1.001126565367507
$ws.Range("F5").Value = 0.995493715611908
$ws.Range("G5").Value = 1.001126565367507
$ws.Range("H5").Value = 1.001126565367507
$ws.Range("I5").Value = 1.003004196338207
$ws.Range("J5").Value = 0.9975215417366059
$ws.Range("K5").Value = 1.001126565367507
$ws.Range("L5").Value = 1.001126565367507
$ws.Range("M5").Value = 0.9983101404897077
$ws.Range("N5").Value = 0.9983101404897077
$ws.Range("O5").Value = 0.9980472742386738
$ws.Range("P5").Value = 0.9992489487823075
$ws.Range("Q5").Value = 0.9992489487823075
$ws.Range("R5").Value = 0.9997183529286076
$ws.Range("S5").Value = 0.9997183529286076
$ws.Range("T5").Value = 0.9998998582982073

# Rename the "Thomas Hex" category to "Matthies Hex" (shared string used at row 11, col B)
$ws.Range("B11").Value = "Matthies Hex"

Write-Host "edit complete"
